$wb = $excel.ActiveWorkbook

$wsAnalysis = $wb.Worksheets.Item("Анализ")
$wsJava = $wb.Worksheets.Item("Java")

# ---------------------------------------------------------------------------
# 1. Update the "Java questions" progress count (Java!A6: 13 -> 14).
#    Everything downstream (Анализ!B3, D3, D4, D5 and the linked chart cache)
#    recalculates automatically off of this cell.
# ---------------------------------------------------------------------------
$wsJava.Range("A6").Value = 14

# ---------------------------------------------------------------------------
# 2. Resize/reposition the six charts on "Анализ" (the right-hand column of
#    charts was narrowed so it no longer overlaps the far columns).
# ---------------------------------------------------------------------------
$chartGeom = @{
    1 = @(323.4291962968, 85.0066929134, 389.0269371309, 132.8098425197)
    2 = @(324.0058104700, 222.1990551181, 390.0404016978, 132.8096062992)
    3 = @(717.4446373647, 85.1777165354, 437.0022576280, 132.8098425197)
    4 = @(718.3098342151, 220.8185039370, 442.9319426673, 132.7499212598)
    5 = @(322.8977789739, 358.5836220472, 390.0002442175, 132.9826771654)
    6 = @(716.3483381521, 358.7029133858, 442.9319426673, 132.9826771654)
}

for ($i = 1; $i -le $wsAnalysis.ChartObjects().Count; $i++) {
    $co = $wsAnalysis.ChartObjects($i)
    $geom = $chartGeom[$i]
    $co.Left = $geom[0]
    $co.Top = $geom[1]
    $co.Width = $geom[2]
    $co.Height = $geom[3]
}

# ---------------------------------------------------------------------------
# 3. Update sheet views: zoom in on "Анализ" and move the selection/active
#    tab from "Анализ" to "Java".
# ---------------------------------------------------------------------------
$wsAnalysis.Activate()
$excel.ActiveWindow.Zoom = 160
$wsAnalysis.Range("O34").Select()

$wsJava.Activate()
$wsJava.Range("F20").Select()
